# aula-08: edit the notes of slide 1 and remove the "Avaliação: Perplexidade" slide.
$p = $ppt.ActivePresentation

# --- 1) Slide 1 notes: the speaker-notes body placeholder gained a run of
#        three tab characters (previously just an empty endParaRPr). ---
$notes1 = $p.Slides.Item(1).NotesPage
$notesBody = $notes1.Shapes.Item(2)
$notesBody.TextFrame.TextRange.Text = "`t`t`t"

# --- 2) Delete the "Avaliação: Perplexidade" slide (slide #36, sldId 299). ---
$p.Slides.Item(36).Delete()
